$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits at the
#    end of the hyperlink paragraph (it will be re-created at the end
#    of the new discussion bullet paragraph instead).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the hyperlink paragraph ("Reference:" link) and the empty
#    trailing paragraph that currently follows it.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$hyperlinkPara = $d.Paragraphs.Item($count - 1)
$trailerPara = $d.Paragraphs.Item($count)

# ------------------------------------------------------------------
# 3. Turn the empty trailing paragraph into the "Discussion:" heading
#    (it already carries the right bold/Verdana formatting).
# ------------------------------------------------------------------
$trailerPara.Range.Text = "Discussion:"

$headingRange = $trailerPara.Range
$headingRange.Font.Name = "Verdana"
$headingRange.Font.Bold = 1

# ------------------------------------------------------------------
# 4. Insert a brand-new paragraph after the heading for the bulleted
#    discussion text.
# ------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.InsertParagraphAfter()

$listPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$listRange = $listPara.Range
$listRange.Text = "Admin will add new employees, their salary and salaryHeads. So, will we store array of in each record for salaryhead of employee?"

# Formatting for the bullet text run: Verdana, not bold, en-US.
$listRange.Font.Name = "Verdana"
$listRange.Font.Bold = 0
$listRange.LanguageID = 1033

# Apply the default round-bullet list (creates a brand new numbered
# list definition, same as picking the Bullet-Library default bullet
# from the Word UI).
$bulletGallery = $word.ListGalleries.Item(1)
$bulletTemplate = $bulletGallery.ListTemplates.Item(1)
$listPara.Range.ListFormat.ApplyListTemplate($bulletTemplate)
$listPara.Style = "List Paragraph"

# ------------------------------------------------------------------
# 5. Re-create the "_GoBack" bookmark, now collapsed at the end of the
#    new discussion bullet paragraph (right before its paragraph mark),
#    matching where it used to sit relative to the previous content.
# ------------------------------------------------------------------
$listPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $listPara2.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)
